$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 109, shifting rows 109..177 down to 110..178,
# preserving all their existing values/formatting.
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new record's data.
$ws.Range("A109").Value = 6
$ws.Range("B109").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C109").Value = "Metropolitana"
$ws.Range("D109").Value = 44606
$ws.Range("E109").Value = 13
$ws.Range("F109").Value = 100112001
$ws.Range("G109").Value = "Berenjena"
$ws.Range("H109").Value = "Sin especificar"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 150
$ws.Range("K109").Value = 10000
$ws.Range("L109").Value = 12000
$ws.Range("M109").Value = 10933
$ws.Range("N109").Value = "$/caja 60 unidades"
$ws.Range("O109").Value = "Región Metropolitana"
$ws.Range("P109").Value = 182
$ws.Range("Q109").Value = 60
$ws.Range("R109").Value = "Hortaliza"
